$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 42
$ws.Range("F5").Value = 72
$ws.Range("F6").Value = 2434
$ws.Range("F7").Value = 49
$ws.Range("F10").Value = 1508
$ws.Range("F12").Value = 602
$ws.Range("F13").Value = 1139
$ws.Range("F14").Value = 1139
$ws.Range("F15").Value = 1156
$ws.Range("F16").Value = 487
$ws.Range("F17").Value = 3469
$ws.Range("F18").Value = 637
$ws.Range("F19").Value = 3260
$ws.Range("F20").Value = 728
$ws.Range("F21").Value = 604
$ws.Range("F22").Value = 21
$ws.Range("F24").Value = 8
$ws.Range("F25").Value = 1102
$ws.Range("F27").Value = 46
$ws.Range("F28").Value = 932
$ws.Range("F29").Value = 909
$ws.Range("F30").Value = 72

$ws = $wb.Worksheets.Item("演出")
$ws.Range("G3").Value = "不可售"
$ws.Range("F8").Value = 82
$ws.Range("F14").Value = 204
$ws.Range("F19").Value = 228
$ws.Range("F20").Value = 167
$ws.Range("F21").Value = 463

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 476

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F8").Value = 42
$ws.Range("G9").Value = "不可售"
$ws.Range("F11").Value = 72
$ws.Range("F12").Value = 476
$ws.Range("F13").Value = 2434
$ws.Range("F15").Value = 49
$ws.Range("F18").Value = 82
$ws.Range("F23").Value = 1508
$ws.Range("F24").Value = 1508
$ws.Range("F27").Value = 1139
$ws.Range("F28").Value = 1139
$ws.Range("F30").Value = 204
$ws.Range("F31").Value = 1156
$ws.Range("F32").Value = 487
$ws.Range("F34").Value = 3469
$ws.Range("F35").Value = 637
$ws.Range("F36").Value = 3260
$ws.Range("F37").Value = 728
$ws.Range("F39").Value = 604
$ws.Range("F41").Value = 1102
$ws.Range("F44").Value = 228
$ws.Range("F45").Value = 167
$ws.Range("F46").Value = 463
$ws.Range("F48").Value = 46
$ws.Range("F49").Value = 932
$ws.Range("F50").Value = 909
$ws.Range("F51").Value = 72
